# BALP 1.1.1 and history file updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.1.0 -> 1.1.1
$ws.Range("B3").Value = "1.1.1"

# Experimental: (empty) -> false
# Assigning the literal text "false" directly would be auto-coerced to a
# Boolean cell by Excel, so instead write it as a formula returning the
# text string "false" and then convert it to a static value in place with
# Copy/PasteSpecial(xlPasteValues) - this yields a genuine text cell.
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null

# Date: 2022-05-04T10:16:52-05:00 -> 2022-10-21T09:04:31-05:00
$ws.Range("B8").Value = "2022-10-21T09:04:31-05:00"
